# Adds a new data row (row 79) to Sheet1, mirroring the structure of the
# existing rows (columns A-J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

# Column D holds a date formatted as plain text (e.g. "2023-08-29"), like the
# other rows in this sheet. Simply assigning the literal string would make
# Excel auto-convert it into a real date serial number, so instead stage the
# text in a scratch cell that is explicitly formatted as Text, copy it, and
# paste only the *value* into the target cell - this keeps D79's own style
# untouched (same default formatting as the rest of row 79) while still
# storing the text verbatim.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2023-08-29"
$scratch.Copy()
$ws.Range("D" + $row).PasteSpecial(-4163, $null, $false, $false)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 0
